# "first unique comment with hash" — add a new row (row 7) to the "哈希"
# (Hash) worksheet describing the "First Unique Character in a String"
# LeetCode problem (#387), matching the table's existing columns:
#   A: No.   B: leetcode id   C: 题目(question)   D: 解题方法(approach)
#   E: 解题关键词(keywords)   F: 时间复杂度(time)   G: 空间复杂度(space)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("哈希")
$ws.Activate()

# --- Row 7 content -------------------------------------------------------
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 387
$ws.Cells.Item(7, 3).Value = "给定一个字符串，找到它的第一个不重复的字符，并返回它的索引。如果不存在，则返回 -1"
$ws.Cells.Item(7, 4).Value = "1 使用哈希表存储每个字符在字符串中的出现次数`n2 遍历每个k-v，如果v是1，就取出相应的字符`n3 判断其在字符串中的位置"
$ws.Cells.Item(7, 5).Value = "哈希表`n字符出现次数"
$ws.Cells.Item(7, 6).Value = "O(N), N是元素个数"

# G7 gets its own look: 微软雅黑 14pt black, left/center aligned, wrapped —
# a new style distinct from the rest of the row (which keeps the sheet's
# usual style).
$g7 = $ws.Cells.Item(7, 7)
$g7.Value = "O(N)"
$g7.Font.Name = "微软雅黑"
$g7.Font.Size = 14
$g7.Font.Color = 0
$g7.HorizontalAlignment = -4131
$g7.VerticalAlignment = -4108
$g7.WrapText = $true

# Row grew tall enough to show the wrapped solution steps in column D.
$ws.Rows.Item(7).RowHeight = 66

# --- View state: leave the cursor on the newly added row -----------------
$ws.Range("D11").Select()
